$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in B9 (this also triggers recalculation of dependent formulas D9/D10)
$ws.Range("B9").Value = 0.467

# Update the active selection to B9 to match the saved view state
$ws.Range("B9").Select()
